# SUCCESS_Advert.docx update per commit "update of job advert"
$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the now-orphaned "_GoBack" bookmark that used to sit right
#    after "Internet of Things (IoT)" / before the closing period, and
#    append the new "26 month fixed term" sentence to that paragraph.
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$rng = $d.Content
$old1 = ")."
$new1 = "). The position is for 26 month fixed term but may be extended by mutual agreement subject to availability of funds."
$rng.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ------------------------------------------------------------------
# 2) Expand the "France and Netherlands" sentence with the partner
#    institutions, and fold in the (reworded) "One of the primary
#    responsibilities..." paragraph that used to follow separately.
# ------------------------------------------------------------------
$old2 = "France and Netherlands.  "
$new2 = "France (Inria, VERIMAG, ENS Paris) and Netherlands (University of Twente).  One of the primary responsibilities of our participation in the project is to validate the scientific and technological innovation: the design, security and privacy analysis, development and validation of an IoT pilot scenario from the healthcare sector. The pilot study is a sensor based monitoring architecture for dementia patients with security critical data and actions. "
$rng2 = $d.Content
$rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# Put the "_GoBack" bookmark back, now after "...security critical"
# (right before " data and actions.") in the rewritten sentence.
$rng3 = $d.Content
$rng3.Find.Execute("security critical", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng3.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng3)

# Remove the now-duplicated lead-in text from the start of the next
# paragraph, leaving it starting with "This position will also include..."
$old3 = "One of the primary responsibilities of our participation in the project is to validate the scientific and technological innovation by pilots: to build and test user-aware security of an IoT pilot scenario from the healthcare sector of a sensor based monitoring architecture for dementia patients with security critical data and actions. "
$rng4 = $d.Content
$rng4.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Closing date: push back from 1 February to 1 March 2017, and
#    make the new date red.
# ------------------------------------------------------------------
$rng5 = $d.Content
$rng5.Find.Execute("Closing date: 1 February 2017", $true, $false, $false, $false, $false, $true, 1, $false, "Closing date: 1 March 2017", 2) | Out-Null

$rng6 = $d.Content
$rng6.Find.Execute("1 March 2017", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng6.Font.Color = 255
